$d = $word.ActiveDocument

# Update the date line
$d.Paragraphs.Item(1).Range.Text = "2025-07-03 Thursday"

# Update the multiplication table cells
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "20×22=440"  # was "80×27=2160"
$t.Cell(1,2).Range.Text = "67×91=6097"  # was "90×46=4140"
$t.Cell(1,3).Range.Text = "20×96=1920"  # was "42×68=2856"
$t.Cell(1,4).Range.Text = "29×25=725"  # was "58×43=2494"
$t.Cell(1,5).Range.Text = "98×14=1372"  # was "39×47=1833"
$t.Cell(5,1).Range.Text = "20×85=1700"  # was "29×26=754"
$t.Cell(5,2).Range.Text = "67×24=1608"  # was "21×48=1008"
$t.Cell(5,3).Range.Text = "15×95=1425"  # was "42×56=2352"
$t.Cell(5,4).Range.Text = "66×11=726"  # was "16×61=976"
$t.Cell(5,5).Range.Text = "47×16=752"  # was "75×26=1950"
$t.Cell(10,1).Range.Text = "28×30=840"  # was "76×43=3268"
$t.Cell(10,2).Range.Text = "21×12=252"  # was "60×33=1980"
$t.Cell(10,3).Range.Text = "69×39=2691"  # was "47×33=1551"
$t.Cell(10,4).Range.Text = "79×57=4503"  # was "61×41=2501"
$t.Cell(10,5).Range.Text = "65×68=4420"  # was "17×33=561"
$t.Cell(15,1).Range.Text = "38×63=2394"  # was "18×83=1494"
$t.Cell(15,2).Range.Text = "36×13=468"  # was "38×41=1558"
$t.Cell(15,3).Range.Text = "33×59=1947"  # was "20×69=1380"
$t.Cell(15,4).Range.Text = "53×67=3551"  # was "80×27=2160"
$t.Cell(15,5).Range.Text = "61×98=5978"  # was "26×26=676"
$t.Cell(20,1).Range.Text = "91×36=3276"  # was "75×59=4425"
$t.Cell(20,2).Range.Text = "56×17=952"  # was "17×34=578"
$t.Cell(20,3).Range.Text = "25×43=1075"  # was "18×95=1710"
$t.Cell(20,4).Range.Text = "32×51=1632"  # was "19×44=836"
$t.Cell(20,5).Range.Text = "14×38=532"  # was "32×19=608"
